$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.2160278745644599
$ws.Range("C2").Value = 0.5226480836236934
$ws.Range("J2").Value = 0.01045296167247387
$ws.Range("P2").Value = 0.1254355400696864
$ws.Range("S2").Value = 0.1254355400696864
$ws.Range("B3").Value = 0.0124223602484472
$ws.Range("C3").Value = 0.04347826086956522
$ws.Range("J3").Value = 0.02484472049689441
$ws.Range("P3").Value = 0.7826086956521739
$ws.Range("S3").Value = 0.1366459627329193
$ws.Range("J4").Value = 0.07843137254901961
$ws.Range("P4").Value = 0.6078431372549019
$ws.Range("S4").Value = 0.3137254901960784
$ws.Range("B6").Value = 0.07924528301886792
$ws.Range("D6").Value = 0.01132075471698113
$ws.Range("F6").Value = 0.0830188679245283
$ws.Range("J6").Value = 0.2
$ws.Range("O6").Value = 0.01886792452830189
$ws.Range("Q6").Value = 0.1169811320754717
$ws.Range("R6").Value = 0.0830188679245283
$ws.Range("S6").Value = 0.4075471698113208
$ws.Range("B7").Value = 0.09777777777777778
$ws.Range("D7").Value = 0.02222222222222222
$ws.Range("F7").Value = 0.06222222222222222
$ws.Range("J7").Value = 0.1422222222222222
$ws.Range("O7").Value = 0.01777777777777778
$ws.Range("Q7").Value = 0.1822222222222222
$ws.Range("R7").Value = 0.1022222222222222
$ws.Range("S7").Value = 0.3733333333333334
$ws.Range("B8").Value = 0.08646616541353383
$ws.Range("D8").Value = 0.02067669172932331
$ws.Range("E8").Value = 0.001879699248120301
$ws.Range("F8").Value = 0.05263157894736842
$ws.Range("J8").Value = 0.1221804511278195
$ws.Range("O8").Value = 0.01691729323308271
$ws.Range("Q8").Value = 0.1672932330827068
$ws.Range("R8").Value = 0.08834586466165413
$ws.Range("S8").Value = 0.443609022556391
$ws.Range("B9").Value = 0.07480314960629922
$ws.Range("D9").Value = 0.01574803149606299
$ws.Range("F9").Value = 0.05118110236220472
$ws.Range("J9").Value = 0.1220472440944882
$ws.Range("O9").Value = 0.01574803149606299
$ws.Range("Q9").Value = 0.1771653543307087
$ws.Range("R9").Value = 0.1141732283464567
$ws.Range("S9").Value = 0.4291338582677166
$ws.Range("B10").Value = 0.09008327024981075
$ws.Range("D10").Value = 0.02195306585919758
$ws.Range("E10").Value = 0.000757002271006813
$ws.Range("F10").Value = 0.08629825889477669
$ws.Range("J10").Value = 0.1158213474640424
$ws.Range("O10").Value = 0.01589704769114307
$ws.Range("Q10").Value = 0.2096896290688872
$ws.Range("R10").Value = 0.07191521574564724
$ws.Range("S10").Value = 0.3875851627554883
$ws.Range("F11").Value = 0.003401360544217687
$ws.Range("G11").Value = 0.1054421768707483
$ws.Range("J11").Value = 0.07482993197278912
$ws.Range("K11").Value = 0.163265306122449
$ws.Range("L11").Value = 0.6394557823129252
$ws.Range("S11").Value = 0.01360544217687075
$ws.Range("G12").Value = 0.7828282828282829
$ws.Range("J12").Value = 0.1363636363636364
$ws.Range("K12").Value = 0.0303030303030303
$ws.Range("L12").Value = 0.0303030303030303
$ws.Range("S12").Value = 0.0202020202020202
$ws.Range("G13").Value = 0.7049180327868853
$ws.Range("J13").Value = 0.2295081967213115
$ws.Range("S13").Value = 0.06557377049180328
$ws.Range("F15").Value = 0.02755905511811024
$ws.Range("H15").Value = 0.1141732283464567
$ws.Range("I15").Value = 0.06299212598425197
$ws.Range("J15").Value = 0.3503937007874016
$ws.Range("K15").Value = 0.07874015748031496
$ws.Range("M15").Value = 0.01181102362204724
$ws.Range("N15").Value = 0.003937007874015748
$ws.Range("O15").Value = 0.09842519685039371
$ws.Range("S15").Value = 0.2519685039370079
$ws.Range("F16").Value = 0.02094240837696335
$ws.Range("H16").Value = 0.1727748691099476
$ws.Range("I16").Value = 0.1361256544502618
$ws.Range("J16").Value = 0.3193717277486911
$ws.Range("K16").Value = 0.07329842931937172
$ws.Range("M16").Value = 0.03664921465968586
$ws.Range("O16").Value = 0.08376963350785341
$ws.Range("S16").Value = 0.1570680628272251
$ws.Range("F17").Value = 0.01659751037344398
$ws.Range("H17").Value = 0.1639004149377593
$ws.Range("I17").Value = 0.1265560165975104
$ws.Range("J17").Value = 0.3796680497925311
$ws.Range("K17").Value = 0.08298755186721991
$ws.Range("M17").Value = 0.01452282157676349
$ws.Range("O17").Value = 0.06639004149377593
$ws.Range("S17").Value = 0.1493775933609958
$ws.Range("F18").Value = 0.01869158878504673
$ws.Range("H18").Value = 0.2009345794392523
$ws.Range("I18").Value = 0.1261682242990654
$ws.Range("J18").Value = 0.3271028037383177
$ws.Range("K18").Value = 0.1121495327102804
$ws.Range("M18").Value = 0.03271028037383177
$ws.Range("O18").Value = 0.05607476635514019
$ws.Range("S18").Value = 0.1261682242990654
$ws.Range("F19").Value = 0.01372549019607843
$ws.Range("H19").Value = 0.2274509803921569
$ws.Range("I19").Value = 0.08235294117647059
$ws.Range("J19").Value = 0.3542483660130719
$ws.Range("K19").Value = 0.09019607843137255
$ws.Range("M19").Value = 0.02549019607843137
$ws.Range("N19").Value = 0.00130718954248366
$ws.Range("O19").Value = 0.06209150326797386
$ws.Range("S19").Value = 0.1431372549019608
